$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before "bitch" (original row 86) with the "bich" entry.
$ws.Rows.Item(86).Insert()
$ws.Range("A86").Value = 330
$ws.Range("B86").Value = "bich"
$ws.Range("C86").Value = 50

# Insert a new row before "itbet" (now shifted down to row 108) with the "it" entry.
$ws.Rows.Item(108).Insert()
$ws.Range("A108").Value = 331
$ws.Range("B108").Value = "it"
$ws.Range("C108").Value = 10
